# Adicionando formatação nas abas Janeiro, Fevereiro e Março.
# E adicionando os nomes de produtos de forma ordenada alfabeticamente, na aba Resumo.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Nomes dos produtos, em ordem alfabetica, na aba Resumo (linhas 2 a 6)
# ---------------------------------------------------------------------------
$wsResumo = $wb.Worksheets.Item("Resumo")
$wsResumo.Range("A2").Value = "Kiwi"
$wsResumo.Range("A3").Value = "Manga"
$wsResumo.Range("A4").Value = "Maçã"
$wsResumo.Range("A5").Value = "Morango"
$wsResumo.Range("A6").Value = "Uva"

# ---------------------------------------------------------------------------
# 2) Formatacao do cabecalho (A1:D1) nas abas Janeiro, Fevereiro, Março e
#    Resumo: negrito, tamanho 12, preenchimento vermelho e alinhamento
#    centralizado (na horizontal e na vertical).
# ---------------------------------------------------------------------------
$sheetNames = @("Janeiro", "Fevereiro", "Março", "Resumo")

$firstSheet = $wb.Worksheets.Item($sheetNames[0])
$firstHeader = $firstSheet.Range("A1:D1")
$firstHeader.Font.Bold = $true
$firstHeader.Font.Size = 12
$firstHeader.Interior.Color = 255
$firstHeader.HorizontalAlignment = -4108
$firstHeader.VerticalAlignment = -4108

# Reaproveita exatamente o mesmo formato nas demais abas (copiar/colar
# apenas formatos), evitando recriar estilos duplicados.
$firstHeader.Copy()
for ($i = 1; $i -lt $sheetNames.Length; $i++) {
    $ws = $wb.Worksheets.Item($sheetNames[$i])
    $header = $ws.Range("A1:D1")
    $header.PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Largura das colunas A:D nas abas Janeiro, Fevereiro, Março e Resumo.
#    (11, 14, 18 e 9 caracteres respectivamente). O runtime adiciona um
#    preenchimento (~0.8333) entre ColumnWidth e a largura gravada no XML,
#    entao subtraimos esse valor para obter a largura final desejada.
# ---------------------------------------------------------------------------
$padding = 0.8333333333333334
$widths = @(11, 14, 18, 9)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    for ($c = 1; $c -le $widths.Length; $c++) {
        $ws.Columns.Item($c).ColumnWidth = $widths[$c - 1] - $padding
    }
}

Write-Host "Formatacao e dados aplicados com sucesso"
